# Mobile Action Plan (MAP) workbook update
# - Mark WSQ Slides 16-19 rows as having a "Slides Folder" (H21:H24)
# - Mark several milestone rows as "Done in project" (H17, H28, H29, H30,
#   H31, H33, H34, H35) - these used to reference the shared "Video" text
# - Two stray cells (H25, H26) get the short text "c"
# - H24 additionally gets an underlined font to flag the signed/packaged app
# - Selection / scroll position moves to reflect where the user was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H17").Value = "Done in project"

$ws.Range("H21").Value = "Slides Folder"
$ws.Range("H22").Value = "Slides Folder"
$ws.Range("H23").Value = "Slides Folder"
$ws.Range("H24").Value = "Slides Folder"
$ws.Range("H24").Font.Underline = $true

$ws.Range("H25").Value = "c"
$ws.Range("H26").Value = "c"

$ws.Range("H28").Value = "Done in project"
$ws.Range("H29").Value = "Done in project"
$ws.Range("H30").Value = "Done in project"
$ws.Range("H31").Value = "Done in project"
$ws.Range("H33").Value = "Done in project"
$ws.Range("H34").Value = "Done in project"
$ws.Range("H35").Value = "Done in project"

# Move the viewport / selection like the author's session ended up
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H24").Select()
